# Scheduled-runner refresh: update cached FFXIV Leve profit figures
# (currentAveragePrice* / LevePrice* / LeveProfit* columns) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to the latest market-board pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 206
$ws.Range("I12").Value = 179.8
$ws.Range("K12").Value = 179.8
$ws.Range("M12").Value = -9.800000000000011
$ws.Range("H40").Value = 4196.4
$ws.Range("I40").Value = 3120.8125
$ws.Range("J40").Value = 6108.5557
$ws.Range("K40").Value = 3120.8125
$ws.Range("L40").Value = 6108.5557
$ws.Range("M40").Value = -2945.8125
$ws.Range("N40").Value = -6458.5557
$ws.Range("H48").Value = 9352.267
$ws.Range("J48").Value = 10420.444
$ws.Range("L48").Value = 31261.332
$ws.Range("N48").Value = -31845.332
$ws.Range("H51").Value = 7412.125
$ws.Range("I51").Value = 6699.3335
$ws.Range("J51").Value = 7839.8
$ws.Range("K51").Value = 6699.3335
$ws.Range("L51").Value = 7839.8
$ws.Range("M51").Value = -6215.3335
$ws.Range("N51").Value = -8807.799999999999
$ws.Range("H56").Value = 9352.267
$ws.Range("J56").Value = 10420.444
$ws.Range("L56").Value = 31261.332
$ws.Range("N56").Value = -32329.332
$ws.Range("H74").Value = 7540
$ws.Range("I74").Value = 4876
$ws.Range("K74").Value = 4876
$ws.Range("M74").Value = -3940
$ws.Range("H77").Value = 7540
$ws.Range("I77").Value = 4876
$ws.Range("K77").Value = 24380
$ws.Range("M77").Value = -19700
$ws.Range("H116").Value = 11409.723
$ws.Range("I116").Value = 10218.5
$ws.Range("J116").Value = 12898.75
$ws.Range("K116").Value = 10218.5
$ws.Range("L116").Value = 12898.75
$ws.Range("M116").Value = -6776.5
$ws.Range("N116").Value = -19782.75
$ws.Range("H132").Value = 1147.8306
$ws.Range("I132").Value = 1161.434
$ws.Range("K132").Value = 3484.302
$ws.Range("M132").Value = -954.3019999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15154309
$ws.Range("I32").Value = 16130990
$ws.Range("K32").Value = 16130990
$ws.Range("M32").Value = -16130703
$ws.Range("H63").Value = 2140.1428
$ws.Range("I63").Value = 1996.8334
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 1996.8334
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -1310.8334
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 2140.1428
$ws.Range("I66").Value = 1996.8334
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 9984.166999999999
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -6552.166999999999
$ws.Range("N66").Value = -21864
$ws.Range("H122").Value = 2127.04
$ws.Range("I122").Value = 2138.3044
$ws.Range("J122").Value = 1997.5
$ws.Range("K122").Value = 6414.9132
$ws.Range("L122").Value = 5992.5
$ws.Range("M122").Value = -3964.9132
$ws.Range("N122").Value = -10892.5
$ws.Range("H132").Value = 2000.3673
$ws.Range("I132").Value = 1648.2727
$ws.Range("K132").Value = 4944.8181
$ws.Range("M132").Value = -2414.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 14433.333
$ws.Range("I82").Value = 3737.875
$ws.Range("J82").Value = 99997
$ws.Range("K82").Value = 3737.875
$ws.Range("L82").Value = 99997
$ws.Range("M82").Value = -3354.875
$ws.Range("N82").Value = -100763
$ws.Range("H85").Value = 14433.333
$ws.Range("I85").Value = 3737.875
$ws.Range("J85").Value = 99997
$ws.Range("K85").Value = 3737.875
$ws.Range("L85").Value = 99997
$ws.Range("M85").Value = -2411.875
$ws.Range("N85").Value = -102649

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 7668.4165
$ws.Range("I105").Value = 7878.8887
$ws.Range("K105").Value = 7878.8887
$ws.Range("M105").Value = -6131.8887
$ws.Range("H132").Value = 3467.3438
$ws.Range("I132").Value = 2351.1155
$ws.Range("K132").Value = 7053.3465
$ws.Range("M132").Value = -4523.3465
$ws.Range("H134").Value = 5122.7085
$ws.Range("I134").Value = 3345.4375
$ws.Range("K134").Value = 10036.3125
$ws.Range("M134").Value = -7501.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 6484517.5
$ws.Range("I34").Value = 8643412
$ws.Range("J34").Value = 7833.3335
$ws.Range("K34").Value = 25930236
$ws.Range("L34").Value = 23500.0005
$ws.Range("M34").Value = -25930152
$ws.Range("N34").Value = -23668.0005
$ws.Range("H36").Value = 20
$ws.Range("I36").Value = 20
$ws.Range("K36").Value = 60
$ws.Range("M36").Value = 109
$ws.Range("H39").Value = 2073.3333
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -15588
$ws.Range("H55").Value = 3066.4
$ws.Range("H125").Value = 4667
$ws.Range("I125").Value = 1500.375
$ws.Range("J125").Value = 30000
$ws.Range("K125").Value = 4501.125
$ws.Range("L125").Value = 90000
$ws.Range("M125").Value = 418.875
$ws.Range("N125").Value = -99840
$ws.Range("H131").Value = 7100577
$ws.Range("J131").Value = 5403037
$ws.Range("L131").Value = 16209111
$ws.Range("N131").Value = -16219191

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1283.3
$ws.Range("J107").Value = 3398.3333
$ws.Range("L107").Value = 3398.3333
$ws.Range("N107").Value = -7238.3333
$ws.Range("H132").Value = 198514.39
$ws.Range("I132").Value = 236722.39
$ws.Range("J132").Value = 4000.9092
$ws.Range("K132").Value = 710167.17
$ws.Range("L132").Value = 12002.7276
$ws.Range("M132").Value = -707637.17
$ws.Range("N132").Value = -17062.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2580.8125
$ws.Range("I46").Value = 733.1667
$ws.Range("J46").Value = 3689.4
$ws.Range("K46").Value = 733.1667
$ws.Range("L46").Value = 3689.4
$ws.Range("M46").Value = -545.1667
$ws.Range("N46").Value = -4065.4
$ws.Range("H122").Value = 5100.846
$ws.Range("I122").Value = 3423
$ws.Range("K122").Value = 10269
$ws.Range("M122").Value = -7819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 22411.2
$ws.Range("I55").Value = 17667.666
$ws.Range("J55").Value = 29526.5
$ws.Range("K55").Value = 17667.666
$ws.Range("L55").Value = 29526.5
$ws.Range("M55").Value = -17390.666
$ws.Range("N55").Value = -30080.5
$ws.Range("H132").Value = 6626.3145
$ws.Range("I132").Value = 2431.1304
$ws.Range("K132").Value = 7293.3912
$ws.Range("M132").Value = -4763.3912
